$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 26831
$ws.Range("E2").Value = 523778156497
$ws.Range("F2").Value = 5423976183
$ws.Range("G2").Value = -0.09397
$ws.Range("D3").Value = 1554.14
$ws.Range("E3").Value = 186950386334
$ws.Range("F3").Value = 2456993720
$ws.Range("G3").Value = 0.48933
$ws.Range("D4").Value = 0.999956
$ws.Range("E4").Value = 83523310091
$ws.Range("F4").Value = 8866877614
$ws.Range("G4").Value = 0.01701
$ws.Range("D5").Value = 206.87
$ws.Range("E5").Value = 31846892603
$ws.Range("F5").Value = 202274670
$ws.Range("G5").Value = 0.51519
$ws.Range("D6").Value = 0.486281
$ws.Range("E6").Value = 25991088961
$ws.Range("F6").Value = 282398047
$ws.Range("G6").Value = 0.12439
$ws.Range("D7").Value = 0.999383
$ws.Range("E7").Value = 25111780254
$ws.Range("F7").Value = 2408120721
$ws.Range("G7").Value = -0.07718
$ws.Range("D8").Value = 1554.33
$ws.Range("E8").Value = 13711685780
$ws.Range("F8").Value = 15178326
$ws.Range("G8").Value = 0.49134
$ws.Range("D9").Value = 21.79
$ws.Range("E9").Value = 9079614066
$ws.Range("F9").Value = 160621207
$ws.Range("G9").Value = -0.46341
$ws.Range("D10").Value = 0.246233
$ws.Range("E10").Value = 8595337336
$ws.Range("F10").Value = 76973323
$ws.Range("G10").Value = 0.00062
$ws.Range("D11").Value = 0.059116
$ws.Range("E11").Value = 8369733123
$ws.Range("F11").Value = 147973120
$ws.Range("G11").Value = 0.45077
$ws.Range("D12").Value = 0.08538
$ws.Range("E12").Value = 7596548689
$ws.Range("F12").Value = 200737055
$ws.Range("G12").Value = 0.03473
$ws.Range("D13").Value = 1.92
$ws.Range("E13").Value = 6628379109
$ws.Range("F13").Value = 13710188
$ws.Range("G13").Value = -1.0927
$ws.Range("B14").Value = "DOT"
$ws.Range("C14").Value = "Polkadot"
$ws.Range("D14").Value = 3.73
$ws.Range("E14").Value = 4798308944
$ws.Range("F14").Value = 67436843
$ws.Range("G14").Value = -0.05501
$ws.Range("B15").Value = "MATIC"
$ws.Range("C15").Value = "Polygon"
$ws.Range("D15").Value = 0.515084
$ws.Range("E15").Value = 4791926070
$ws.Range("F15").Value = 93527363
$ws.Range("G15").Value = -0.3185
$ws.Range("D16").Value = 61.46
$ws.Range("E16").Value = 4534635047
$ws.Range("F16").Value = 145957665
$ws.Range("G16").Value = -0.04898
$ws.Range("D17").Value = 26817
$ws.Range("E17").Value = 4371497770
$ws.Range("F17").Value = 27278860
$ws.Range("G17").Value = -0.17271
$ws.Range("D18").Value = 213.45
$ws.Range("E18").Value = 4170310656
$ws.Range("F18").Value = 67161082
$ws.Range("G18").Value = -1.15309
$ws.Range("D19").Value = 0.00000699
$ws.Range("E19").Value = 4120307764
$ws.Range("F19").Value = 73670495
$ws.Range("G19").Value = -1.85399
$ws.Range("D20").Value = 7.34
$ws.Range("E20").Value = 4089650602
$ws.Range("F20").Value = 108370812
$ws.Range("G20").Value = 1.35632
$ws.Range("D21").Value = 0.999561
$ws.Range("E21").Value = 3783747143
$ws.Range("F21").Value = 40242195
$ws.Range("G21").Value = -0.01461
$ws.Range("D22").Value = 3.81
$ws.Range("E22").Value = 3536976356
$ws.Range("F22").Value = 171107
$ws.Range("G22").Value = 2.12704
$ws.Range("B23").Value = "TUSD"
$ws.Range("C23").Value = "TrueUSD"
$ws.Range("D23").Value = 0.999133
$ws.Range("E23").Value = 3366138724
$ws.Range("F23").Value = 56196315
$ws.Range("G23").Value = -0.00339
$ws.Range("B24").Value = "AVAX"
$ws.Range("C24").Value = "Avalanche"
$ws.Range("D24").Value = 9.1
$ws.Range("E24").Value = 3231314749
$ws.Range("F24").Value = 79655091
$ws.Range("G24").Value = -0.82112
$ws.Range("D25").Value = 4.11
$ws.Range("E25").Value = 3097965772
$ws.Range("F25").Value = 36833040
$ws.Range("G25").Value = 0.81024
$ws.Range("D26").Value = 0.104614
$ws.Range("E26").Value = 2906201000
$ws.Range("F26").Value = 37239204
$ws.Range("G26").Value = 0.8154400000000001
$ws.Range("D27").Value = 152.63
$ws.Range("E27").Value = 2772332911
$ws.Range("F27").Value = 33984604
$ws.Range("G27").Value = -0.21978
$ws.Range("D28").Value = 43.06
$ws.Range("E28").Value = 2584193964
$ws.Range("F28").Value = 2063446
$ws.Range("G28").Value = -0.40202
$ws.Range("B29").Value = "ETC"
$ws.Range("C29").Value = "Ethereum Classic"
$ws.Range("D29").Value = 14.9
$ws.Range("E29").Value = 2134682284
$ws.Range("F29").Value = 37499873
$ws.Range("G29").Value = 0.14464
$ws.Range("B30").Value = "BUSD"
$ws.Range("C30").Value = "BUSD"
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 2124771813
$ws.Range("F30").Value = 951097462
$ws.Range("G30").Value = -0.01297
$ws.Range("D31").Value = 6.55
$ws.Range("E31").Value = 1917984094
$ws.Range("F31").Value = 64769227
$ws.Range("G31").Value = -0.91562
$ws.Range("D32").Value = 0.04690354
$ws.Range("E32").Value = 1570965561
$ws.Range("F32").Value = 19818919
$ws.Range("G32").Value = 0.26215
$ws.Range("D33").Value = 3.2
$ws.Range("E33").Value = 1468035047
$ws.Range("F33").Value = 34954945
$ws.Range("G33").Value = -0.08821
$ws.Range("D34").Value = 1.6
$ws.Range("E34").Value = 1421785067
$ws.Range("F34").Value = 14852466
$ws.Range("G34").Value = 0.45377
$ws.Range("D35").Value = 3.16
$ws.Range("E35").Value = 1415363775
$ws.Range("F35").Value = 24730898
$ws.Range("G35").Value = 2.09883
$ws.Range("D36").Value = 0.051391
$ws.Range("E36").Value = 1353443301
$ws.Range("F36").Value = 3633620
$ws.Range("G36").Value = -0.59581
$ws.Range("B37").Value = "MKR"
$ws.Range("C37").Value = "Maker"
$ws.Range("D37").Value = 1417.72
$ws.Range("E37").Value = 1276913760
$ws.Range("F37").Value = 37041590
$ws.Range("G37").Value = 0.91169
$ws.Range("B38").Value = "QNT"
$ws.Range("C38").Value = "Quant"
$ws.Range("D38").Value = 86.12
$ws.Range("E38").Value = 1252838514
$ws.Range("F38").Value = 9424440
$ws.Range("G38").Value = -0.86638
$ws.Range("B39").Value = "VET"
$ws.Range("C39").Value = "VeChain"
$ws.Range("D39").Value = 0.01656644
$ws.Range("E39").Value = 1204534632
$ws.Range("F39").Value = 24668204
$ws.Range("G39").Value = 0.5394
$ws.Range("B40").Value = "APT"
$ws.Range("C40").Value = "Aptos"
$ws.Range("D40").Value = 4.89
$ws.Range("E40").Value = 1195886411
$ws.Range("F40").Value = 29687707
$ws.Range("G40").Value = -0.6807299999999999
$ws.Range("B41").Value = "OP"
$ws.Range("C41").Value = "Optimism"
$ws.Range("D41").Value = 1.2
$ws.Range("E41").Value = 1055268557
$ws.Range("F41").Value = 37201743
$ws.Range("G41").Value = -0.41282
$ws.Range("B42").Value = "ARB"
$ws.Range("C42").Value = "Arbitrum"
$ws.Range("D42").Value = 0.803496
$ws.Range("E42").Value = 1024802318
$ws.Range("F42").Value = 45444086
$ws.Range("G42").Value = -0.16054
$ws.Range("B43").Value = "NEAR"
$ws.Range("C43").Value = "NEAR Protocol"
$ws.Range("D43").Value = 1.017
$ws.Range("E43").Value = 1002174436
$ws.Range("F43").Value = 28109989
$ws.Range("G43").Value = 0.58691
$ws.Range("B44").Value = "MNT"
$ws.Range("C44").Value = "Mantle"
$ws.Range("D44").Value = 0.322835
$ws.Range("E44").Value = 1001508467
$ws.Range("F44").Value = 26859380
$ws.Range("G44").Value = -0.89456
$ws.Range("B45").Value = "AAVE"
$ws.Range("C45").Value = "Aave"
$ws.Range("D45").Value = 63.94
$ws.Range("E45").Value = 932703645
$ws.Range("F45").Value = 39716994
$ws.Range("G45").Value = -0.44775
$ws.Range("B46").Value = "RETH"
$ws.Range("C46").Value = "Rocket Pool ETH"
$ws.Range("D46").Value = 1691.63
$ws.Range("E46").Value = 907461020
$ws.Range("F46").Value = 4276457
$ws.Range("G46").Value = 0.55769
$ws.Range("B47").Value = "KAS"
$ws.Range("C47").Value = "Kaspa"
$ws.Range("D47").Value = 0.04282639
$ws.Range("E47").Value = 903409493
$ws.Range("F47").Value = 6053272
$ws.Range("G47").Value = -0.34977
$ws.Range("B48").Value = "GRT"
$ws.Range("C48").Value = "The Graph"
$ws.Range("D48").Value = 0.082245
$ws.Range("E48").Value = 760500385
$ws.Range("F48").Value = 15582143
$ws.Range("G48").Value = 1.16705
$ws.Range("B49").Value = "ALGO"
$ws.Range("C49").Value = "Algorand"
$ws.Range("D49").Value = 0.09535100000000001
$ws.Range("E49").Value = 756286461
$ws.Range("F49").Value = 18430589
$ws.Range("G49").Value = 0.09649000000000001
$ws.Range("D50").Value = 5.09
$ws.Range("E50").Value = 733969660
$ws.Range("F50").Value = 9486386
$ws.Range("G50").Value = 1.56067
$ws.Range("D51").Value = 1.001
$ws.Range("E51").Value = 726907273
$ws.Range("F51").Value = 42838119
$ws.Range("G51").Value = 0.1704
